$wb = $excel.ActiveWorkbook
$picarro = $wb.Worksheets.Item("Picarro")
$lgr = $wb.Worksheets.Item("LGR")
$newSheet = $wb.Worksheets.Add($null, $lgr)
$newSheet.Name = "Test1"
$picarro.Range("A1:F13").Copy($newSheet.Range("A1")) | Out-Null
